$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 91, shifting rows 91-119 down to 92-120
$ws.Rows.Item(91).Insert()

# Populate the new row 91 with data
$ws.Range("A91").Value = 1
$ws.Range("B91").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C91").Value = "Arica y Parinacota"
$ws.Range("D91").Value = 44798
$ws.Range("D91").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E91").Value = 15
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100106
$ws.Range("H91").Value = "Oleaginosos"
$ws.Range("I91").Value = 100106002
$ws.Range("J91").Value = "Palta"
$ws.Range("K91").Value = "Hass"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 400
$ws.Range("N91").Value = 24000
$ws.Range("O91").Value = 25000
$ws.Range("P91").Value = 24500
$ws.Range("Q91").Value = "$/bandeja 10 kilos"
$ws.Range("R91").Value = "Perú"
$ws.Range("S91").Value = 2450
$ws.Range("T91").Value = 10
